$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RespRate (column D) values for rows 2-5 from 27 to 8.
# Downstream formulas in A, B, C, G recalc automatically.
$ws.Range("D2:D5").Value = 8

# Move the active selection to D6, matching the saved view state.
$ws.Range("D6").Select()
